$wb = $excel.ActiveWorkbook

# Sheet names affected by this data refresh: "展览" and "全部类型"
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 1891
    $ws.Range("F4").Value = 1158
    $ws.Range("F5").Value = 1191
    $ws.Range("F7").Value = 5990
}
